$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 98, shifting existing rows 98-106 down to 99-107.
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new weekly price record.
$ws.Range("A98").Value = 9
$ws.Range("B98").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C98").Value = "Metropolitana"
$ws.Range("D98").Value = 44610
$ws.Range("E98").Value = 13
$ws.Range("F98").Value = 100112022
$ws.Range("G98").Value = "Arveja Verde"
$ws.Range("H98").Value = "Sin especificar"
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 27
$ws.Range("K98").Value = 23000
$ws.Range("L98").Value = 25000
$ws.Range("M98").Value = 23963
$ws.Range("N98").Value = "`$/saco 25 kilos"
$ws.Range("O98").Value = "Carahue"
$ws.Range("P98").Value = 959
$ws.Range("Q98").Value = 25
$ws.Range("R98").Value = "Hortaliza"

# Preserve the date number format on the inserted row's date cell,
# matching the format used by the other rows in column D (s="2").
$ws.Range("D98").NumberFormat = $ws.Range("D99").NumberFormat
